$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO... ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911773017745.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911773267658.csv"
$ws1.Range("B4").Value = "go_stims-16502911773277676.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911773437765.csv"

# --- Sheet 2: NB_TO... ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16502911791717668.csv"
$ws2.Range("B3").Value = "ZB-match_0-16502911777737684.csv"
$ws2.Range("B4").Value = "TB-16502911834116452.csv"
$ws2.Range("B5").Value = "ZB-match_6-16502911775067704.csv"
$ws2.Range("B6").Value = "TB-16502911824367692.csv"
$ws2.Range("B7").Value = "TB-16502911819667675.csv"
$ws2.Range("B8").Value = "OB-16502911801957664.csv"
$ws2.Range("B9").Value = "ZB-match_7-1650291177844772.csv"
$ws2.Range("B10").Value = "OB-16502911795717752.csv"

# --- Sheet 3: RS_TO... ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO... ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502911834706507.csv"
$ws4.Range("B3").Value = "ZM_stims-1650291183444652.csv"
$ws4.Range("B4").Value = "MM_stims-16502911835026646.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911834716532.csv"
$ws4.Range("B6").Value = "MM_stims-16502911835176597.csv"
$ws4.Range("B7").Value = "ZM_stims-16502911835046518.csv"

# --- Sheet 5: vSAT_TO... ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16502911835496495.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502911835646522.csv"
$ws5.Range("B4").Value = "SAT_stims-1650291183524656.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502911835786524.csv"

# --- Rename sheets (must happen after referencing by index, names are changing) ---
$ws1.Name = "GNG_TO-16502911773447845"
$ws2.Name = "NB_TO-16502911834386508"
$ws3.Name = "RS_TO-16502911834416542"
$ws4.Name = "TOL_TO-16502911835196517"
$ws5.Name = "vSAT_TO-16502911835946493"
